# Started processing factor data into modeling data
# - Drop the 5th "Developed" quartile bucket (column F) and the 5th
#   "Emerging" quartile bucket (former column K), collapsing each
#   Market group from 5 columns down to 4.
# - Refresh the forward-return values in row 4 with the newly
#   recomputed figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra quartile column from each 5-wide group (Developed: B:F,
# Emerging: G:K). Deleting column F first collapses Developed to B:E and
# shifts Emerging left to F:J; deleting the new last column (J, originally K)
# then collapses Emerging to F:I. Excel automatically keeps the dimension,
# row spans and merged header cells (B1:E1 / F1:I1) in sync.
$ws.Columns("F:F").Delete()
$ws.Columns("J:J").Delete()

# Update the forward-return (row 4) figures with the newly processed values.
$ws.Range("B4").Value = 0.009706056663553892
$ws.Range("C4").Value = 0.006442343349931703
$ws.Range("D4").Value = 0.007185107563260074
$ws.Range("E4").Value = 0.007281605013398909
$ws.Range("F4").Value = 0.01407919400705496
$ws.Range("G4").Value = 0.00907865458304232
$ws.Range("H4").Value = 0.01022902647427289
$ws.Range("I4").Value = 0.01316605673043128
